$wb = $excel.ActiveWorkbook

function Set-TextValue($range, $text) {
    # Force the cell to be written as literal text (shared string) rather
    # than letting Excel auto-detect date-looking strings and convert them
    # into date serial numbers. Resetting the style back to Normal afterwards
    # keeps the cell's formatting identical to the surrounding un-styled cells.
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.Style = "Normal"
}

# --- Sheet "Stato Attuale": update operator assignments ---
$stato = $wb.Worksheets.Item("Stato Attuale")

# Row 13: targa GL342AD - operator changes from DARIO.LELLA to
# "LELLA DARIO (INCIDENTATA FARE PRATICA - DECUNCIA CC)"; new assignment date set.
$stato.Range("B13").Value = "LELLA DARIO (INCIDENTATA FARE PRATICA - DECUNCIA CC)"
Set-TextValue $stato.Range("C13") "2026-01-29"

# Row 42: targa GL350TJ - operator changes from DI IANNI FELICIA to LELLA DARIO;
# assignment date updated to the new change date.
$stato.Range("B42").Value = "LELLA DARIO"
Set-TextValue $stato.Range("C42") "2026-01-29"

# --- Sheet "Storico Passaggi": log the two operator changes ---
$storico = $wb.Worksheets.Item("Storico Passaggi")

$storico.Range("A2").Value = "GL350TJ"
$storico.Range("B2").Value = "DI IANNI FELICIA"
$storico.Range("C2").Value = "LELLA DARIO"
Set-TextValue $storico.Range("D2") "2026-01-29"

$storico.Range("A3").Value = "GL342AD"
$storico.Range("B3").Value = "DARIO.LELLA"
$storico.Range("C3").Value = "LELLA DARIO (INCIDENTATA FARE PRATICA - DECUNCIA CC)"
Set-TextValue $storico.Range("D3") "2026-01-29"
